$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 holds a sequence of "type: blog" cards (one per day-grid column).
# A brand-new blog card (ser: 145) is being slotted in at column C, which
# bumps every later blog card one slot to the right:
#   C10 (ser: 140) -> G10
#   G10 (ser: 144) -> J10
#   J10 (ser: 143) -> dropped off the end
# C10 itself becomes the new "ser: 145" card.

$oldC10 = $ws.Range("C10").Value2
$oldG10 = $ws.Range("G10").Value2

$ws.Range("J10").Value = $oldG10
$ws.Range("G10").Value = $oldC10
$ws.Range("C10").Value = "type: blog`nwidth: 2`nheight: 1`nser: 145"

# Match the saved selection state (cell M10 active).
$ws.Range("M10").Select()
